# ICTU-Kwaliteitsaanpak-Checklist.xlsx edit
# Commit: "Add missing Scrum terms. Fixed #90."
#
# Summary of change:
#  - Bump the version/date mentioned in the intro cell (A1).
#  - Rework the bullet list describing the mandatory Scrum ingredients
#    (M05 section): add "sprint retrospective" to the process bullet,
#    merge "Definition of Ready"/"Definition of Done" into a single
#    bullet, merge "Product backlog"/"sprint backlog" into a single
#    bullet - which means one row (the old "5. Product backlog" row)
#    disappears from the table entirely.
#  - Update the M05 cell-comment text to match.
#  - Delete the now-redundant row and let everything below shift up by one
#    row (mergeCells / dataValidation / conditionalFormatting / dimension
#    all need to reflect the new, smaller range).
#  - Cell comments are anchored to fixed cell refs in this engine (they do
#    NOT automatically follow a Rows.Delete shift), so every comment that
#    lived below the deleted row has to be recreated one row higher.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the version / date banner in A1.
$ws.Range("A1").Value = "Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de Kwaliteitsaanpak ICTU Software Realisatie versie 1.2.1-build.1, 29-08-2018."

# 2) Rework the four (was five) Scrum bullet cells, B31:B35.
#    B31 (item 1) is unchanged.
$ws.Range("B32").Value = "2. Proces met daily scrum, sprints, sprint planning, sprint review, sprint retrospective en sprint refinement"
$ws.Range("B33").Value = "3. Definition of Ready en Definition of Done"
$ws.Range("B34").Value = "4. Product backlog en sprint backlog"

# 3) Capture every comment that sits on or below row 36 - it needs to be
#    recreated one row higher once row 35 is removed. Capture first, since
#    once the comments are deleted their text would be lost.
$commentRefsOld = @("B36","B37","B46","B47","B48","B49","B50","B52","B53","B54","B55","B56","B57","B67","B71","B72","B73","B74","B75")
$commentTexts = @{}
foreach ($ref in $commentRefsOld) {
    $cmt = $ws.Range($ref).Comment
    $commentTexts[$ref] = $cmt.Text()
}

# Also grab/update the B30 (M05) comment text - the bullet list inside it
# needs the same rewording as the table rows.
$b30Old = $ws.Range("B30").Comment.Text()
$oldBullets = "- Scrum team bestaand uit product owner, ontwikkelteam en Scrum master,`n- Proces met daily scrum, sprints, sprint planning, sprint review, sprint refinement,`n- Definition of Done,`n- Definition of Ready,`n- Product backlog."
$newBullets = "- Scrum team bestaand uit product owner, ontwikkelteam en Scrum master,`n- Proces met daily scrum, sprints, sprint planning, sprint review, sprint retrospective en sprint refinement,`n- Definition of Ready en Definition of Done,`n- Product backlog en sprint backlog."
$b30New = $b30Old.Replace($oldBullets, $newBullets)

# 4) Delete the now-redundant row (old "5. Product backlog" row); everything
#    below shifts up by one.
$ws.Rows(35).Delete()

# 5) Re-apply the updated M05 comment text (row 30 itself didn't move).
$ws.Range("B30").Comment.Text($b30New)

# 6) Recreate every captured comment one row higher than where it used to be.
foreach ($ref in $commentRefsOld) {
    $oldRowCol = $ref -replace '[0-9]+$',''
    $oldRowNum = [int]($ref -replace '^[A-Z]+','')
    $newRef = "$oldRowCol$($oldRowNum - 1)"
    $target = $ws.Range($newRef)
    if ($target.Comment) {
        $target.Comment.Delete()
    }
    $target.AddComment($commentTexts[$ref]) | Out-Null
}
